$d = $word.ActiveDocument

$d.Content.Find.Execute("74×23=", $true, $false, $false, $false, $false, $true, 1, $false, "37×96=", 2)
$d.Content.Find.Execute("91×68=", $true, $false, $false, $false, $false, $true, 1, $false, "38×59=", 2)
$d.Content.Find.Execute("92×43=", $true, $false, $false, $false, $false, $true, 1, $false, "12×38=", 2)
$d.Content.Find.Execute("68×43=", $true, $false, $false, $false, $false, $true, 1, $false, "85×90=", 2)
$d.Content.Find.Execute("11×33=", $true, $false, $false, $false, $false, $true, 1, $false, "55×12=", 2)
$d.Content.Find.Execute("53×65=", $true, $false, $false, $false, $false, $true, 1, $false, "58×92=", 2)
$d.Content.Find.Execute("21×90=", $true, $false, $false, $false, $false, $true, 1, $false, "40×65=", 2)
$d.Content.Find.Execute("95×77=", $true, $false, $false, $false, $false, $true, 1, $false, "20×93=", 2)
$d.Content.Find.Execute("30×18=", $true, $false, $false, $false, $false, $true, 1, $false, "98×38=", 2)
$d.Content.Find.Execute("48×27=", $true, $false, $false, $false, $false, $true, 1, $false, "25×49=", 2)
$d.Content.Find.Execute("48×67=", $true, $false, $false, $false, $false, $true, 1, $false, "36×74=", 2)
$d.Content.Find.Execute("29×98=", $true, $false, $false, $false, $false, $true, 1, $false, "24×24=", 2)
$d.Content.Find.Execute("24×98=", $true, $false, $false, $false, $false, $true, 1, $false, "53×50=", 2)
$d.Content.Find.Execute("90×36=", $true, $false, $false, $false, $false, $true, 1, $false, "81×37=", 2)
$d.Content.Find.Execute("58×86=", $true, $false, $false, $false, $false, $true, 1, $false, "90×47=", 2)
$d.Content.Find.Execute("55×70=", $true, $false, $false, $false, $false, $true, 1, $false, "95×65=", 2)
$d.Content.Find.Execute("23×45=", $true, $false, $false, $false, $false, $true, 1, $false, "30×71=", 2)
$d.Content.Find.Execute("79×17=", $true, $false, $false, $false, $false, $true, 1, $false, "24×25=", 2)
$d.Content.Find.Execute("40×42=", $true, $false, $false, $false, $false, $true, 1, $false, "94×86=", 2)
$d.Content.Find.Execute("20×99=", $true, $false, $false, $false, $false, $true, 1, $false, "80×66=", 2)
$d.Content.Find.Execute("93×24=", $true, $false, $false, $false, $false, $true, 1, $false, "84×75=", 2)
$d.Content.Find.Execute("99×76=", $true, $false, $false, $false, $false, $true, 1, $false, "77×49=", 2)
$d.Content.Find.Execute("73×59=", $true, $false, $false, $false, $false, $true, 1, $false, "13×61=", 2)
$d.Content.Find.Execute("64×37=", $true, $false, $false, $false, $false, $true, 1, $false, "52×66=", 2)
$d.Content.Find.Execute("33×11=", $true, $false, $false, $false, $false, $true, 1, $false, "82×87=", 2)
